# Teaching_assignments.xlsx - add new Programming-language entries to column B
# (Scala, Erlang, Clojure, Go, R, Rust, Elixir) and move the active selection
# to B13, reflecting the repository name change / directory restructuring
# that introduced these new skill rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "Scala"
$ws.Range("B8").Value  = "Erlang"
$ws.Range("B9").Value  = "Clojure"
$ws.Range("B10").Value = "Go"
$ws.Range("B11").Value = "R"
$ws.Range("B12").Value = "Rust"
$ws.Range("B13").Value = "Elixir (On Erlang VM)"

[void]$ws.Range("B13").Select()
